$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = '27.089.79'
$ws.Range("E2").Value = '  -1.09%  '

$ws.Range("D3").Value = '1.823.68'
$ws.Range("E3").Value = '  -1.00%  '

$ws.Range("E4").Value = '  -0.45%  '

$ws.Range("D5").Value = '311.23'
$ws.Range("E5").Value = '  -1.29%  '

$ws.Range("E6").Value = '  -0.49%  '

$ws.Range("D7").Value = '0.4656'
$ws.Range("E7").Value = '  -1.89%  '

$ws.Range("D8").Value = '0.3638'
$ws.Range("E8").Value = '  -1.79%  '

$ws.Range("D9").Value = '0.07295'
$ws.Range("E9").Value = '  -2.37%  '

$ws.Range("D10").Value = '0.8693'
$ws.Range("E10").Value = '  -1.92%  '

$ws.Range("D11").Value = '20.21'
$ws.Range("E11").Value = '  -1.54%  '

$ws.Range("D12").Value = '0.07618'
$ws.Range("E12").Value = '  +3.41%  '

$ws.Range("D13").Value = '1.845.74'
$ws.Range("E13").Value = '  -0.70%  '

$ws.Range("D14").Value = '92.93'
$ws.Range("E14").Value = '  -0.44%  '

$ws.Range("D15").Value = '5.341'
$ws.Range("E15").Value = '  -2.70%  '

$ws.Range("D16").Value = '6.477'
$ws.Range("E16").Value = '  -1.58%  '

$ws.Range("E17").Value = '  -0.58%  '

$ws.Range("D18").Value = '0.000008645'
$ws.Range("E18").Value = '  -2.51%  '

$ws.Range("E19").Value = '  -0.44%  '

$ws.Range("D20").Value = '27.271.04'
$ws.Range("E20").Value = '  -0.52%  '

$ws.Range("D21").Value = '14.49'
$ws.Range("E21").Value = '  -2.55%  '

$ws.Range("D22").Value = '5.192'
$ws.Range("E22").Value = '  -3.08%  '

$ws.Range("E23").Value = '  -1.56%  '

$ws.Range("D24").Value = '2.076.26'
$ws.Range("E24").Value = '  -0.14%  '

$ws.Range("D25").Value = '151.70'
$ws.Range("E25").Value = '  -0.25%  '

$ws.Range("D26").Value = '1.866'
$ws.Range("E26").Value = '  -2.31%  '

$ws.Range("D27").Value = '18.20'
$ws.Range("E27").Value = '  -2.51%  '

$ws.Range("D28").Value = '2.113'
$ws.Range("E28").Value = '  -3.25%  '

$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").Value = '5.095'
$ws.Range("E29").Value = '  -3.57%  '

$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").Value = '115.94'
$ws.Range("E30").Value = '  -1.88%  '

$ws.Range("D31").Value = '0.08923'
$ws.Range("E31").Value = '  -0.51%  '

$ws.Range("D32").Value = '2.956'
$ws.Range("E32").Value = '  +0.30%  '

$ws.Range("D33").Value = '0.7348'
$ws.Range("E33").Value = '  -3.66%  '

$ws.Range("E34").Value = '  -2.55%  '

$ws.Range("D35").Value = '1.141'
$ws.Range("E35").Value = '  -3.32%  '

$ws.Range("D36").Value = '1.008'
$ws.Range("E36").Value = '  -0.48%  '

$ws.Range("D37").Value = '2.556'
$ws.Range("E37").Value = '  +7.18%  '

$ws.Range("D38").Value = '0.05262'
$ws.Range("E38").Value = '  -2.11%  '

$ws.Range("D39").Value = '1.070'
$ws.Range("E39").Value = '  -3.15%  '

$ws.Range("D40").Value = '0.01917'
$ws.Range("E40").Value = '  -2.35%  '

$ws.Range("D41").Value = '2.936'
$ws.Range("E41").Value = '  -2.06%  '

$ws.Range("D42").Value = '7.131'
$ws.Range("E42").Value = '  -2.64%  '

$ws.Range("D43").Value = '0.5218'
$ws.Range("E43").Value = '  -2.69%  '

$ws.Range("D44").Value = '0.1632'
$ws.Range("E44").Value = '  -2.22%  '

$ws.Range("D45").Value = '8.272'
$ws.Range("E45").Value = '  -3.38%  '

$ws.Range("D46").Value = '0.4873'
$ws.Range("E46").Value = '  -2.27%  '

$ws.Range("D47").Value = '1.008'
$ws.Range("E47").Value = '  -0.59%  '

$ws.Range("D48").Value = '103.79'
$ws.Range("E48").Value = '  -1.33%  '

$ws.Range("D49").Value = '10.14'
$ws.Range("E49").Value = '  -3.57%  '

$ws.Range("D50").Value = '1.637'
$ws.Range("E50").Value = '  -2.68%  '

$ws.Range("D51").Value = '0.06246'
$ws.Range("E51").Value = '  -1.34%  '
